$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1442
    $ws.Range("F5").Value = 76
    $ws.Range("F6").Value = 16
}
